$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet1"

# Rewrite headers (row 1) so shared strings are created/ordered as:
# VENDA, CLIENTE, PAGAMENTO, VALOR, ID_PROD, QNT, VALOR_VENDA
$ws.Range("A1").Value = "VENDA"
$ws.Range("B1").Value = "CLIENTE"
$ws.Range("C1").Value = "PAGAMENTO"
$ws.Range("D1").Value = "VALOR"
$ws.Range("E1").Value = "ID_PROD"
$ws.Range("F1").Value = "QNT"
$ws.Range("G1").Value = "VALOR_VENDA"

# Row 2: first sale (single product)
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 55323285
# Leading apostrophe -> quote-prefixed text cell (shared string "a_vista")
$ws.Cells.Item(2,3).Value = "'a_vista"
# Quantity stored as a number but displayed with a Text number format
$ws.Cells.Item(2,4).Value = 3
$ws.Cells.Item(2,4).NumberFormat = "@"

# Row 3: second order, first product line
$ws.Cells.Item(3,5).Value = 88645377
$ws.Cells.Item(3,5).NumberFormat = "General"
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1

# Row 4: second order, second product line (same order id, new product)
$ws.Cells.Item(4,5).Value = 88645401
$ws.Cells.Item(4,5).NumberFormat = "General"
$ws.Cells.Item(4,6).Value = 2
$ws.Cells.Item(4,7).Value = 1

# Widen the new ID_PROD-order column (G)
$ws.Columns.Item(7).ColumnWidth = 15.6

# Restore selection similar to the authored workbook
$ws.Range("E14").Select() | Out-Null
